$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "219.61").
# Force the cell format to Text before assigning so Excel keeps it as a
# string (matching the original inline-string cell) instead of silently
# converting it to a floating point number.
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "26.493.96"
$ws.Range("E2").Value = "  +1.80%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.671.78"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("E4").Value = "  +0.07%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "219.61"
$ws.Range("E5").Value = "  +2.11%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.5287"
$ws.Range("E6").Value = "  +1.19%  "
$ws.Range("E7").Value = "  +0.08%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06375"
$ws.Range("E9").Value = "  +0.27%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "21.77"
$ws.Range("E10").Value = "  +4.70%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07797"
$ws.Range("E11").Value = "  +1.61%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.672.17"
$ws.Range("E12").Value = "  +1.60%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "4.480"
$ws.Range("E13").Value = "  +1.30%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.5574"
$ws.Range("E14").Value = "  +0.62%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.0₅8303"
$ws.Range("E15").Value = "  -0.36%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "65.55"
$ws.Range("E16").Value = "  +1.17%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "26.506.89"
$ws.Range("E17").Value = "  +1.83%  "
$ws.Range("E18").Value = "  +0.05%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "4.760"
$ws.Range("E19").Value = "  +0.97%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "193.21"
$ws.Range("E20").Value = "  +2.64%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "10.33"
$ws.Range("E21").Value = "  +1.64%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.303"
$ws.Range("E22").Value = "  +0.70%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "1.003"
$ws.Range("E23").Value = "  +0.12%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.1268"
$ws.Range("E24").Value = "  +3.93%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "138.44"
$ws.Range("E25").Value = "  -4.35%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "7.387"
$ws.Range("E26").Value = "  -0.20%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "16.31"
$ws.Range("E27").Value = "  +3.07%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.425"
$ws.Range("E28").Value = "  +2.68%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.06223"
$ws.Range("E29").Value = "  +4.45%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.288"
$ws.Range("E30").Value = "  +1.84%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.607"
$ws.Range("E31").Value = "  +5.98%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.419"
$ws.Range("E32").Value = "  +0.65%  "
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("E34").Value = "  +1.22%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.6122"
$ws.Range("E35").Value = "  +8.86%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.414"
$ws.Range("E36").Value = "  +0.93%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.782"
$ws.Range("E37").Value = "  +1.01%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.01616"
$ws.Range("E38").Value = "  +0.59%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "6.045"
$ws.Range("E39").Value = "  +3.24%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.091.95"
$ws.Range("E40").Value = "  +6.34%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("E42").Value = "  +0.02%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "100.56"
$ws.Range("E43").Value = "  +1.94%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.816.70"
$ws.Range("E44").Value = "  +1.23%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "58.30"
$ws.Range("E45").Value = "  +4.70%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "8.137"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.520"
$ws.Range("E48").Value = "  +10.31%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.05194"
$ws.Range("E49").Value = "  +1.00%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "6.005"
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.4231"
$ws.Range("E51").Value = "  +0.38%  "
